# Elimina EC anteriores y se agregan nuevos, se modifica base de datos
# Rewrite the "Periodo Mora" (E) and "Valor Mora" (F) columns of the account
# statement table (rows 16-28) with the updated periods/values, keeping the
# worker identity columns (B/C/D) aligned to the new grouping described by
# the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New table contents for rows 16-28 (columns B:G)
# B = Tipo Doc, C = N Doc Trabajador, D = Nombre Trabajador,
# E = Periodo Mora, F = Valor Mora, G = Salario Basico
$rows = @(
    @{ Row = 16; B = "CC"; C = "45539915";   D = "MAYLIN JOHANNA VERBEL AYOLA";    E = "2310"; F = 35574; G = 1160000 },
    @{ Row = 17; B = "CC"; C = "45539915";   D = "MAYLIN JOHANNA VERBEL AYOLA";    E = "2309"; F = 46400; G = 1160000 },
    @{ Row = 18; B = "CC"; C = "45539915";   D = "MAYLIN JOHANNA VERBEL AYOLA";    E = "2308"; F = 46400; G = 1160000 },
    @{ Row = 19; B = "CC"; C = "45539915";   D = "MAYLIN JOHANNA VERBEL AYOLA";    E = "2307"; F = 46400; G = 1160000 },
    @{ Row = 20; B = "CC"; C = "45539915";   D = "MAYLIN JOHANNA VERBEL AYOLA";    E = "2306"; F = 46400; G = 1160000 },
    @{ Row = 21; B = "CC"; C = "45539915";   D = "MAYLIN JOHANNA VERBEL AYOLA";    E = "2305"; F = 46400; G = 1160000 },
    @{ Row = 22; B = "CC"; C = "45539915";   D = "MAYLIN JOHANNA VERBEL AYOLA";    E = "2304"; F = 46400; G = 1160000 },
    @{ Row = 23; B = "CC"; C = "1002198845"; D = "WILSON DAVID VILLADIEGO VERBEL"; E = "2310"; F = 35574; G = 1160000 },
    @{ Row = 24; B = "CC"; C = "1002198845"; D = "WILSON DAVID VILLADIEGO VERBEL"; E = "2309"; F = 46400; G = 1160000 },
    @{ Row = 25; B = "CC"; C = "1002198845"; D = "WILSON DAVID VILLADIEGO VERBEL"; E = "2307"; F = 46400; G = 1160000 },
    @{ Row = 26; B = "CC"; C = "1002198845"; D = "WILSON DAVID VILLADIEGO VERBEL"; E = "2306"; F = 46400; G = 1160000 },
    @{ Row = 27; B = "CC"; C = "1002198845"; D = "WILSON DAVID VILLADIEGO VERBEL"; E = "2305"; F = 46400; G = 1160000 },
    @{ Row = 28; B = "CC"; C = "1002198845"; D = "WILSON DAVID VILLADIEGO VERBEL"; E = "2304"; F = 46400; G = 1160000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
    $ws.Range("G$n").Value = $r.G
}
